$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content updates -------------------------------------------------

# B14: "Alternativa 1\n[Dados inválidos]\n(Passo 5)" -> "Exceção 1\n[Dados inválidos]\n(Passo 5)"
$ws.Range("B14").Value = "Exceção 1" + [char]10 + "[Dados inválidos]" + [char]10 + "(Passo 5)"

# D14: "5.1. Informa que login falhou" -> "5.1. Informa que as credenciais são inválidas falhou"
$ws.Range("D14").Value = "5.1. Informa que as credenciais são inválidas falhou"

# D15: "Regressa a 4" -> cleared
$ws.Range("D15").Value = ""

# --- Column width -----------------------------------------------------------
# Target stored width is 50.125 chars; this engine quantizes ColumnWidth to
# 1/6-character steps, so 49.333333 is the closest input that lands on the
# nearest achievable stored width (50.1666...).
$ws.Columns.Item(4).ColumnWidth = 49.333333333333336

# --- New rows 17-19 ----------------------------------------------------------

$ws.Rows.Item(17).RowHeight = 20.25
$ws.Rows.Item(18).RowHeight = 18.75
$ws.Rows.Item(19).RowHeight = 21.75

# Formatting: copy existing equivalent styles then strip the border so we end
# up with fillId=2 (light gray header shading) + borderId=0 variants, the
# same combinations used by the new cellXfs entries in the target workbook.

$ws.Range("B6").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B17").Borders.LineStyle = -4142

$ws.Range("B10").Copy()
$ws.Range("B18:B19").PasteSpecial(-4122)
$ws.Range("B18:B19").Borders.LineStyle = -4142

$ws.Range("A1").Copy()
$ws.Range("C17:E19").PasteSpecial(-4122)
$ws.Range("C17:E19").Borders.LineStyle = -4142

$ws.Application.CutCopyMode = $false

# Merge the new left-hand label cell (mirrors B14:B16 above it)
$ws.Range("B17:B19").Merge()

# --- Selection ---------------------------------------------------------------
$ws.Range("C18").Select()
